$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-16 in place with refreshed TPM-derived NATMI statistics.
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnc"
$ws.Cells.Item(2, 3).Value = "Ptprz1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 11.5183295
$ws.Cells.Item(2, 8).Value = 23.036659
$ws.Cells.Item(2, 9).Value = 0.06427717328589268
$ws.Cells.Item(2, 10).Value = 0.05102762964987022
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 0.071111
$ws.Cells.Item(2, 14).Value = 0.142222
$ws.Cells.Item(2, 15).Value = 0.02711460746047303
$ws.Cells.Item(2, 16).Value = 0.02678527708115022
$ws.Cells.Item(2, 17).Value = 0.8190799290744999
$ws.Cells.Item(2, 18).Value = 3.276319716298
$ws.Cells.Item(2, 19).Value = 0.001742850322315784
$ws.Cells.Item(2, 20).Value = 0.00136678919896609

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnc"
$ws.Cells.Item(3, 3).Value = "Ptprz1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 11.5183295
$ws.Cells.Item(3, 8).Value = 23.036659
$ws.Cells.Item(3, 9).Value = 0.06427717328589268
$ws.Cells.Item(3, 10).Value = 0.05102762964987022
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.064491
$ws.Cells.Item(3, 14).Value = 0.193473
$ws.Cells.Item(3, 15).Value = 0.0245904030281302
$ws.Cells.Item(3, 16).Value = 0.03643759694506741
$ws.Cells.Item(3, 17).Value = 0.7428285877845001
$ws.Cells.Item(3, 18).Value = 4.456971526707
$ws.Cells.Item(3, 19).Value = 0.001580601596609065
$ws.Cells.Item(3, 20).Value = 0.001859324202244142

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnc"
$ws.Cells.Item(4, 3).Value = "Ptprz1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 11.5183295
$ws.Cells.Item(4, 8).Value = 23.036659
$ws.Cells.Item(4, 9).Value = 0.06427717328589268
$ws.Cells.Item(4, 10).Value = 0.05102762964987022
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.4870065
$ws.Cells.Item(4, 14).Value = 4.974013
$ws.Cells.Item(4, 15).Value = 0.9482949895113968
$ws.Cells.Item(4, 16).Value = 0.9367771259737823
$ws.Cells.Item(4, 17).Value = 28.64616033564175
$ws.Cells.Item(4, 18).Value = 114.584641342567
$ws.Cells.Item(4, 19).Value = 0.06095372136696783
$ws.Cells.Item(4, 20).Value = 0.04780151624865998

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnc"
$ws.Cells.Item(5, 3).Value = "Ptprz1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 93.03757333333333
$ws.Cells.Item(5, 8).Value = 279.11272
$ws.Cells.Item(5, 9).Value = 0.5191891952080042
$ws.Cells.Item(5, 10).Value = 0.6182520002891011
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.5
$ws.Cells.Item(5, 13).Value = 0.071111
$ws.Cells.Item(5, 14).Value = 0.142222
$ws.Cells.Item(5, 15).Value = 0.02711460746047303
$ws.Cells.Item(5, 16).Value = 0.02678527708115022
$ws.Cells.Item(5, 17).Value = 6.615994877306665
$ws.Cells.Item(5, 18).Value = 39.69596926383999
$ws.Cells.Item(5, 19).Value = 0.01407761122578394
$ws.Cells.Item(5, 20).Value = 0.01656005113371894

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnc"
$ws.Cells.Item(6, 3).Value = "Ptprz1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 93.03757333333333
$ws.Cells.Item(6, 8).Value = 279.11272
$ws.Cells.Item(6, 9).Value = 0.5191891952080042
$ws.Cells.Item(6, 10).Value = 0.6182520002891011
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.064491
$ws.Cells.Item(6, 14).Value = 0.193473
$ws.Cells.Item(6, 15).Value = 0.0245904030281302
$ws.Cells.Item(6, 16).Value = 0.03643759694506741
$ws.Cells.Item(6, 17).Value = 6.000086141840001
$ws.Cells.Item(6, 18).Value = 54.00077527656
$ws.Cells.Item(6, 19).Value = 0.01276707155801539
$ws.Cells.Item(6, 20).Value = 0.02252761719701597

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnc"
$ws.Cells.Item(7, 3).Value = "Ptprz1"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 93.03757333333333
$ws.Cells.Item(7, 8).Value = 279.11272
$ws.Cells.Item(7, 9).Value = 0.5191891952080042
$ws.Cells.Item(7, 10).Value = 0.6182520002891011
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.4870065
$ws.Cells.Item(7, 14).Value = 4.974013
$ws.Cells.Item(7, 15).Value = 0.9482949895113968
$ws.Cells.Item(7, 16).Value = 0.9367771259737823
$ws.Cells.Item(7, 17).Value = 231.3850496242267
$ws.Cells.Item(7, 18).Value = 1388.31029774536
$ws.Cells.Item(7, 19).Value = 0.4923445124242049
$ws.Cells.Item(7, 20).Value = 0.5791643319583661

$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Tnc"
$ws.Cells.Item(8, 3).Value = "Ptprz1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 74.6205215
$ws.Cells.Item(8, 8).Value = 149.241043
$ws.Cells.Item(8, 9).Value = 0.4164142197129523
$ws.Cells.Item(8, 10).Value = 0.3305781741512238
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.5
$ws.Cells.Item(8, 13).Value = 0.071111
$ws.Cells.Item(8, 14).Value = 0.142222
$ws.Cells.Item(8, 15).Value = 0.02711460746047303
$ws.Cells.Item(8, 16).Value = 0.02678527708115022
$ws.Cells.Item(8, 17).Value = 5.306339904386499
$ws.Cells.Item(8, 18).Value = 21.225359617546
$ws.Cells.Item(8, 19).Value = 0.01129090810847587
$ws.Cells.Item(8, 20).Value = 0.008854627991621261

$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Tnc"
$ws.Cells.Item(9, 3).Value = "Ptprz1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 74.6205215
$ws.Cells.Item(9, 8).Value = 149.241043
$ws.Cells.Item(9, 9).Value = 0.4164142197129523
$ws.Cells.Item(9, 10).Value = 0.3305781741512238
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.064491
$ws.Cells.Item(9, 14).Value = 0.193473
$ws.Cells.Item(9, 15).Value = 0.0245904030281302
$ws.Cells.Item(9, 16).Value = 0.03643759694506741
$ws.Cells.Item(9, 17).Value = 4.8123520520565
$ws.Cells.Item(9, 18).Value = 28.874112312339
$ws.Cells.Item(9, 19).Value = 0.01023979348938586
$ws.Cells.Item(9, 20).Value = 0.01204547426855859

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Tnc"
$ws.Cells.Item(10, 3).Value = "Ptprz1"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 74.6205215
$ws.Cells.Item(10, 8).Value = 149.241043
$ws.Cells.Item(10, 9).Value = 0.4164142197129523
$ws.Cells.Item(10, 10).Value = 0.3305781741512238
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.4870065
$ws.Cells.Item(10, 14).Value = 4.974013
$ws.Cells.Item(10, 15).Value = 0.9482949895113968
$ws.Cells.Item(10, 16).Value = 0.9367771259737823
$ws.Cells.Item(10, 17).Value = 185.5817220038898
$ws.Cells.Item(10, 18).Value = 742.326888015559
$ws.Cells.Item(10, 19).Value = 0.3948835181150905
$ws.Cells.Item(10, 20).Value = 0.3096780718910438

$ws.Cells.Item(11, 1).Value = "Neutrophils"
$ws.Cells.Item(11, 2).Value = "Tnc"
$ws.Cells.Item(11, 3).Value = "Ptprz1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.008067
$ws.Cells.Item(11, 8).Value = 0.024201
$ws.Cells.Item(11, 9).Value = 0.00004501728804487631
$ws.Cells.Item(11, 10).Value = 0.00005360671723953153
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.5
$ws.Cells.Item(11, 13).Value = 0.071111
$ws.Cells.Item(11, 14).Value = 0.142222
$ws.Cells.Item(11, 15).Value = 0.02711460746047303
$ws.Cells.Item(11, 16).Value = 0.02678527708115022
$ws.Cells.Item(11, 17).Value = 0.0005736524369999999
$ws.Cells.Item(11, 18).Value = 0.003441914622
$ws.Cells.Item(11, 19).Value = 0.000001220626094271867
$ws.Cells.Item(11, 20).Value = 0.000001435870774671724

$ws.Cells.Item(12, 1).Value = "Neutrophils"
$ws.Cells.Item(12, 2).Value = "Tnc"
$ws.Cells.Item(12, 3).Value = "Ptprz1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.008067
$ws.Cells.Item(12, 8).Value = 0.024201
$ws.Cells.Item(12, 9).Value = 0.00004501728804487631
$ws.Cells.Item(12, 10).Value = 0.00005360671723953153
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.064491
$ws.Cells.Item(12, 14).Value = 0.193473
$ws.Cells.Item(12, 15).Value = 0.0245904030281302
$ws.Cells.Item(12, 16).Value = 0.03643759694506741
$ws.Cells.Item(12, 17).Value = 0.0005202488970000001
$ws.Cells.Item(12, 18).Value = 0.004682240073
$ws.Cells.Item(12, 19).Value = 0.000001106993256256936
$ws.Cells.Item(12, 20).Value = 0.000001953299956322246

$ws.Cells.Item(13, 1).Value = "Neutrophils"
$ws.Cells.Item(13, 2).Value = "Tnc"
$ws.Cells.Item(13, 3).Value = "Ptprz1"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.008067
$ws.Cells.Item(13, 8).Value = 0.024201
$ws.Cells.Item(13, 9).Value = 0.00004501728804487631
$ws.Cells.Item(13, 10).Value = 0.00005360671723953153
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.4870065
$ws.Cells.Item(13, 14).Value = 4.974013
$ws.Cells.Item(13, 15).Value = 0.9482949895113968
$ws.Cells.Item(13, 16).Value = 0.9367771259737823
$ws.Cells.Item(13, 17).Value = 0.0200626814355
$ws.Cells.Item(13, 18).Value = 0.120376088613
$ws.Cells.Item(13, 19).Value = 0.00004268966869434751
$ws.Cells.Item(13, 20).Value = 0.00005021754650853755

$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Tnc"
$ws.Cells.Item(14, 3).Value = "Ptprz1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.01333133333333333
$ws.Cells.Item(14, 8).Value = 0.039994
$ws.Cells.Item(14, 9).Value = 0.00007439450510585445
$ws.Cells.Item(14, 10).Value = 0.00008858919256550655
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.5
$ws.Cells.Item(14, 13).Value = 0.071111
$ws.Cells.Item(14, 14).Value = 0.142222
$ws.Cells.Item(14, 15).Value = 0.02711460746047303
$ws.Cells.Item(14, 16).Value = 0.02678527708115022
$ws.Cells.Item(14, 17).Value = 0.0009480044446666667
$ws.Cells.Item(14, 18).Value = 0.005688026668
$ws.Cells.Item(14, 19).Value = 0.0000020171778031614
$ws.Cells.Item(14, 20).Value = 0.000002372886069262466

$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Tnc"
$ws.Cells.Item(15, 3).Value = "Ptprz1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.01333133333333333
$ws.Cells.Item(15, 8).Value = 0.039994
$ws.Cells.Item(15, 9).Value = 0.00007439450510585445
$ws.Cells.Item(15, 10).Value = 0.00008858919256550655
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.064491
$ws.Cells.Item(15, 14).Value = 0.193473
$ws.Cells.Item(15, 15).Value = 0.0245904030281302
$ws.Cells.Item(15, 16).Value = 0.03643759694506741
$ws.Cells.Item(15, 17).Value = 0.0008597510180000001
$ws.Cells.Item(15, 18).Value = 0.007737759162000001
$ws.Cells.Item(15, 19).Value = 0.000001829390863631251
$ws.Cells.Item(15, 20).Value = 0.00000322797729239089

$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Tnc"
$ws.Cells.Item(16, 3).Value = "Ptprz1"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.01333133333333333
$ws.Cells.Item(16, 8).Value = 0.039994
$ws.Cells.Item(16, 9).Value = 0.00007439450510585445
$ws.Cells.Item(16, 10).Value = 0.00008858919256550655
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 2.4870065
$ws.Cells.Item(16, 14).Value = 4.974013
$ws.Cells.Item(16, 15).Value = 0.9482949895113968
$ws.Cells.Item(16, 16).Value = 0.9367771259737823
$ws.Cells.Item(16, 17).Value = 0.03315511265366667
$ws.Cells.Item(16, 18).Value = 0.198930675922
$ws.Cells.Item(16, 19).Value = 0.00007054793643906179
$ws.Cells.Item(16, 20).Value = 0.00008298832920385318

# Rows 17-19 (the "Resolving-Mac" block) are dropped entirely; shift rows up and
# shrink the used range back down to A1:T16.
$ws.Range("A17:T19").Delete() | Out-Null
